$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Most cells are plain text/percent strings (with spaces) or numeric-looking
# strings whose General-format number display matches the literal text,
# so a direct .Value assignment reproduces the same stored text.

$ws.Range("D2").Value = '27.531.51'
$ws.Range("E2").Value = '  +0.52%  '
$ws.Range("D3").Value = '1.638.37'
$ws.Range("E3").Value = '  -0.78%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").Value = '212.65'
$ws.Range("E5").Value = '  -0.43%  '
$ws.Range("D6").Value = '0.534'
$ws.Range("E6").Value = '  +4.98%  '
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("D8").Value = '22.92'
$ws.Range("E8").Value = '  -5.27%  '
$ws.Range("E9").Value = '  -1.86%  '
$ws.Range("E10").Value = '  -0.76%  '
# D11: force text so a trailing zero ("0.0890") is not dropped by numeric parsing
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0890'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.39%  '
$ws.Range("D12").Value = '1.870.49'
$ws.Range("E12").Value = '  -0.76%  '
$ws.Range("D13").Value = '1.643.85'
$ws.Range("E13").Value = '  -0.10%  '
$ws.Range("E14").Value = '  -1.42%  '
$ws.Range("E15").Value = '  -2.00%  '
$ws.Range("D16").Value = '64.27'
$ws.Range("E16").Value = '  -2.44%  '
$ws.Range("D17").Value = '27.474.98'
$ws.Range("E17").Value = '  +0.30%  '
$ws.Range("D18").Value = '229.63'
$ws.Range("E18").Value = '  -1.81%  '
# D19: force text so a trailing zero ("7.70") is not dropped by numeric parsing
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.70'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.35%  '
$ws.Range("D20").Value = '0.0₃0723'
$ws.Range("E20").Value = '  -0.61%  '
$ws.Range("E21").Value = '  +0.13%  '
# D22: force text so a trailing zero ("4.30") is not dropped by numeric parsing
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.30'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.00%  '
# D23: force text so a trailing zero ("9.90") is not dropped by numeric parsing
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.90'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +6.27%  '
$ws.Range("D24").Value = '1.95'
$ws.Range("E24").Value = '  -3.37%  '
$ws.Range("D25").Value = '149.71'
$ws.Range("E25").Value = '  +1.92%  '
$ws.Range("D26").Value = '6.96'
$ws.Range("D27").Value = '0.112'
$ws.Range("E27").Value = '  +1.41%  '
$ws.Range("E28").Value = '  +0.09%  '
$ws.Range("E29").Value = '  -2.67%  '
$ws.Range("D30").Value = '1.19'
$ws.Range("E30").Value = '  -0.82%  '
$ws.Range("E31").Value = '  -1.97%  '
$ws.Range("D32").Value = '3.28'
$ws.Range("E32").Value = '  -0.56%  '
$ws.Range("D33").Value = '3.16'
$ws.Range("E33").Value = '  +2.08%  '
$ws.Range("D34").Value = '1.422.33'
$ws.Range("E34").Value = '  -2.64%  '
$ws.Range("E35").Value = '  +2.27%  '
$ws.Range("D36").Value = '2.35'
$ws.Range("E36").Value = '  -1.58%  '
$ws.Range("E37").Value = '  -0.37%  '
$ws.Range("D38").Value = '0.875'
$ws.Range("E38").Value = '  -3.71%  '
$ws.Range("E39").Value = '  -2.08%  '
$ws.Range("D40").Value = '0.883'
$ws.Range("E40").Value = '  +12.47%  '
$ws.Range("D41").Value = '1.03'
$ws.Range("E41").Value = '  -0.88%  '
$ws.Range("E42").Value = '  +0.07%  '
$ws.Range("E43").Value = '  -0.38%  '
$ws.Range("D44").Value = '5.52'
$ws.Range("E44").Value = '  +0.96%  '
$ws.Range("E45").Value = '  +1.43%  '
$ws.Range("D46").Value = '64.83'
$ws.Range("E46").Value = '  -1.09%  '
$ws.Range("D47").Value = '1.779.93'
$ws.Range("E47").Value = '  -0.76%  '
$ws.Range("E48").Value = '  -3.41%  '
$ws.Range("D49").Value = '86.13'
$ws.Range("E49").Value = '  -2.63%  '
$ws.Range("D50").Value = '0.0₆0106'
$ws.Range("E50").Value = '  +0.01%  '
$ws.Range("D51").Value = '0.0987'
$ws.Range("E51").Value = '  -2.43%  '
